# Weekly update: insert a new weekly record for "Ajo" (Macroferia Regional
# de Talca) right after row 285, shifting the existing rows 286-299 down to
# 287-300, and populate the newly inserted row 286 with the new week's data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 286; this shifts rows 286..299 down to 287..300
# and grows the used range to row 300.
$ws.Rows.Item(286).Insert()

# Fill in the new row 286 with this week's record.
$ws.Cells.Item(286, 1).Value2  = 5
$ws.Cells.Item(286, 2).Value2  = "Macroferia Regional de Talca"
$ws.Cells.Item(286, 3).Value2  = "Maule"
$ws.Cells.Item(286, 4).Value2  = 44706
$ws.Cells.Item(286, 5).Value2  = 7
$ws.Cells.Item(286, 6).Value2  = 100112003
$ws.Cells.Item(286, 7).Value2  = "Ajo"
$ws.Cells.Item(286, 8).Value2  = "Chino"
$ws.Cells.Item(286, 9).Value2  = "Primera"
$ws.Cells.Item(286, 10).Value2 = 300
$ws.Cells.Item(286, 11).Value2 = 18000
$ws.Cells.Item(286, 12).Value2 = 18000
$ws.Cells.Item(286, 13).Value2 = 18000
$ws.Cells.Item(286, 14).Value2 = "$/caja 10 kilos"
$ws.Cells.Item(286, 15).Value2 = "China"
$ws.Cells.Item(286, 16).Value2 = 1800
$ws.Cells.Item(286, 17).Value2 = 10
$ws.Cells.Item(286, 18).Value2 = "Hortaliza"
